$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B5").Value = 683
$ws.Range("B6").Value = 468
$ws.Range("B16").Value = 437
$ws.Range("B17").Value = 274
$ws.Range("B18").Value = 482
$ws.Range("B19").Value = 410
$ws.Range("B20").Value = 426
$ws.Range("B21").Value = 350
$ws.Range("B22").Value = 398
$ws.Range("B23").Value = 389
$ws.Range("B24").Value = 383
$ws.Range("B25").Value = 303
$ws.Range("B26").Value = 377
$ws.Range("B28").Value = 255
$ws.Range("B29").Value = 350
$ws.Range("B30").Value = 358
$ws.Range("B39").Value = 514
$ws.Range("B50").Value = 427
